$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column D that splits out / normalizes the professor's full
# name (stored in uppercase in column B: "APELLIDOS, NOMBRE") into a
# proper-cased version using the PROPER() worksheet function.
#
# The fill was done in two passes (D1 alone, then D2:D65, then D66:D89),
# which is why the saved workbook ends up with two separate shared-formula
# groups - replicate that here so the resulting formula layout matches.

$ws.Range("D1").Formula = "=PROPER(B1)"
$ws.Range("D2:D65").Formula = "=PROPER(B2)"
$ws.Range("D66:D89").Formula = "=PROPER(B66)"

# Restore the selection to match the state left after performing the fill.
$ws.Range("D1:D89").Select()
